$wb = $excel.ActiveWorkbook

# Each of the 4 sheets ("MID_LFT_#1", "MID_LFT_#2", "MID_PLT_#1", "MID_PLT_#2")
# gets one new daily record appended as row 99 (the sheets previously ended
# at row 98). Copy the last row (98) down to row 99 first so the new row
# inherits the same formatting/style (e.g. the date style on column A),
# then overwrite the cells whose values differ for the new day.

# Sheet 1: MID_LFT_#1
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A98:I98").Copy($ws1.Range("A99:I99"))
$ws1.Cells.Item(99, 1).Value = 45885.46673611111

# Sheet 2: MID_LFT_#2
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A98:I98").Copy($ws2.Range("A99:I99"))
$ws2.Cells.Item(99, 1).Value = 45885.46673611111
$ws2.Cells.Item(99, 4).Value = "0x01,0x20"
$ws2.Cells.Item(99, 8).Value = 288

# Sheet 3: MID_PLT_#1
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A98:I98").Copy($ws3.Range("A99:I99"))
$ws3.Cells.Item(99, 1).Value = 45885.46673611111

# Sheet 4: MID_PLT_#2
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("A98:I98").Copy($ws4.Range("A99:I99"))
$ws4.Cells.Item(99, 1).Value = 45885.46673611111
$ws4.Cells.Item(99, 4).Value = "0x00,0x73"
$ws4.Cells.Item(99, 8).Value = 115
